$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O6").Value = 0
$ws.Range("O7").Value = 0.0009999275207519531
$ws.Range("O10").Value = 0
$ws.Range("O11").Value = 0.01605916023254395
$ws.Range("O12").Value = 0.02594470977783203
$ws.Range("O13").Value = 0.005932807922363281
$ws.Range("O14").Value = 0.0774693489074707
$ws.Range("O17").Value = 0.04897141456604004
$ws.Range("O18").Value = 0
$ws.Range("O19").Value = 0.002053976058959961
$ws.Range("O20").Value = 0.02978205680847168
$ws.Range("O21").Value = 0.006841659545898438
$ws.Range("O23").Value = 0.01694416999816895
$ws.Range("O25").Value = 0.005001544952392578
$ws.Range("O26").Value = 0.0167233943939209
$ws.Range("O27").Value = 0
$ws.Range("O29").Value = 0.05733132362365723
$ws.Range("O30").Value = 0.0213463306427002
$ws.Range("O31").Value = 0.02427220344543457
$ws.Range("O32").Value = 0.01630258560180664
$ws.Range("O33").Value = 0.0009999275207519531
$ws.Range("O34").Value = 0.01202583312988281
$ws.Range("O35").Value = 0
$ws.Range("O36").Value = 0.002997636795043945
$ws.Range("O37").Value = 0.02444839477539062
$ws.Range("O38").Value = 0.01499128341674805
$ws.Range("O41").Value = 0.01192498207092285
$ws.Range("O44").Value = 0
$ws.Range("O48").Value = 0.01692819595336914
$ws.Range("O49").Value = 0.009581327438354492
$ws.Range("O50").Value = 0.00203394889831543
$ws.Range("O53").Value = 0.01462578773498535
$ws.Range("O54").Value = 0.0009992122650146484
$ws.Range("O57").Value = 0.003154993057250977
$ws.Range("O58").Value = 0.0008606910705566406
$ws.Range("O59").Value = 0
$ws.Range("O60").Value = 0.005997180938720703
$ws.Range("O61").Value = 0.000682830810546875
$ws.Range("O62").Value = 0.002043008804321289
$ws.Range("O63").Value = 0.001001119613647461
$ws.Range("O64").Value = 0.01990151405334473
$ws.Range("O65").Value = 0.1624987125396729
$ws.Range("O66").Value = 0.03593325614929199
$ws.Range("O67").Value = 0.003111124038696289
$ws.Range("O68").Value = 0.03015756607055664
$ws.Range("O69").Value = 0
$ws.Range("O70").Value = 0.01012277603149414
$ws.Range("O71").Value = 1.693661689758301
$ws.Range("O72").Value = 3.157593727111816
$ws.Range("O73").Value = 0.04002022743225098
$ws.Range("O75").Value = 0
$ws.Range("O76").Value = 0.0009951591491699219
$ws.Range("O78").Value = 0.02655601501464844
$ws.Range("O79").Value = 0.01642394065856934
$ws.Range("O80").Value = 0.003007173538208008
$ws.Range("O81").Value = 0.2855050563812256
$ws.Range("O82").Value = 0.04378747940063477
$ws.Range("O83").Value = 0.1323230266571045
$ws.Range("O84").Value = 13.45798468589783
$ws.Range("O85").Value = 0.1207764148712158
$ws.Range("O86").Value = 0.01582813262939453
$ws.Range("O88").Value = 0.0159459114074707
$ws.Range("O90").Value = 0.1421191692352295
$ws.Range("O91").Value = 0.008006572723388672
$ws.Range("O92").Value = 0.1319520473480225
$ws.Range("O94").Value = 0
$ws.Range("O95").Value = 0.001003026962280273
$ws.Range("O96").Value = 0.07441973686218262
$ws.Range("O97").Value = 0
$ws.Range("O98").Value = 0.002004861831665039
$ws.Range("O99").Value = 0.001997709274291992
$ws.Range("O100").Value = 0.05442190170288086
$ws.Range("O103").Value = 0.0009992122650146484
$ws.Range("O104").Value = 0.02299046516418457
$ws.Range("O105").Value = 0
$ws.Range("O106").Value = 0.01633620262145996
$ws.Range("O107").Value = 0.03007650375366211
$ws.Range("O108").Value = 0.02606344223022461
$ws.Range("O109").Value = 0.001001358032226562
$ws.Range("O110").Value = 0.001001596450805664
$ws.Range("O111").Value = 0.002002954483032227
$ws.Range("O112").Value = 0.001991033554077148
$ws.Range("O113").Value = 0.006363391876220703
$ws.Range("O114").Value = 0.001998662948608398
$ws.Range("O115").Value = 0.003018617630004883
$ws.Range("O116").Value = 0.03776764869689941
$ws.Range("O117").Value = 0.001088857650756836
$ws.Range("O118").Value = 0.0005028247833251953
$ws.Range("O119").Value = 0
$ws.Range("O120").Value = 0.00100398063659668
$ws.Range("O122").Value = 0.0145106315612793
$ws.Range("O123").Value = 0
$ws.Range("O126").Value = 0
$ws.Range("O127").Value = 0.0005195140838623047
$ws.Range("O128").Value = 0.0009992122650146484
$ws.Range("O129").Value = 0.0009965896606445312
$ws.Range("O130").Value = 0
$ws.Range("O132").Value = 0.0007159709930419922
$ws.Range("O133").Value = 0
$ws.Range("O136").Value = 0.001001358032226562
$ws.Range("O140").Value = 0.001006126403808594
$ws.Range("O141").Value = 0.001611709594726562
$ws.Range("O143").Value = 0
$ws.Range("O146").Value = 0.0008375644683837891
$ws.Range("O147").Value = 0
$ws.Range("O148").Value = 0.00602269172668457
$ws.Range("O149").Value = 0.01870250701904297
$ws.Range("O150").Value = 0.02051353454589844
$ws.Range("O151").Value = 0.05555319786071777
$ws.Range("O152").Value = 0.0009152889251708984
$ws.Range("O153").Value = 0
$ws.Range("O155").Value = 0.0009999275207519531
$ws.Range("O159").Value = 0.04638528823852539
$ws.Range("O160").Value = 0
$ws.Range("O161").Value = 0.008521795272827148
$ws.Range("O163").Value = 0.001996040344238281
$ws.Range("O164").Value = 0.0009965896606445312
$ws.Range("O166").Value = 0.002916574478149414
$ws.Range("O167").Value = 0.001841306686401367
$ws.Range("O168").Value = 0
$ws.Range("O169").Value = 0.009624958038330078
$ws.Range("O170").Value = 0.01565456390380859
$ws.Range("O171").Value = 0.04068517684936523
$ws.Range("O172").Value = 0.009683609008789062
$ws.Range("O173").Value = 0
$ws.Range("O174").Value = 0.0101017951965332
$ws.Range("O175").Value = 0
$ws.Range("O176").Value = 0
$ws.Range("O177").Value = 0
